$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# Step 1: remove the paragraph "This PDF version is provided under the same
# license." entirely (its content+mark disappear, merging into the
# preceding paragraph which keeps its own paragraph mark).
# ---------------------------------------------------------------------------
$pdfIdx = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    if ($d.Paragraphs.Item($i).Range.Text -like "This PDF version*") { $pdfIdx = $i; break }
}
if ($pdfIdx -ne $null) {
    $d.Paragraphs.Item($pdfIdx).Range.Delete()
}

# ---------------------------------------------------------------------------
# Step 2: remove the "License Information" Heading2 paragraph entirely.
# ---------------------------------------------------------------------------
$licIdx = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    if ($d.Paragraphs.Item($i).Range.Text -like "License Information*") { $licIdx = $i; break }
}
if ($licIdx -ne $null) {
    $d.Paragraphs.Item($licIdx).Range.Delete()
}

# ---------------------------------------------------------------------------
# Step 3: rewrite the resource-description paragraph (the one that used to
# start with the bold resource title and continue "... is based on: ...").
# ---------------------------------------------------------------------------
$basedIdx = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    if ($d.Paragraphs.Item($i).Range.Text -like "*is based on*") { $basedIdx = $i; break }
}
$p = $d.Paragraphs.Item($basedIdx)
$pStart = $p.Range.Start
$pEnd = $p.Range.End

# Find the bold run's character span inside the paragraph (it directly
# follows a zero-width empty run, so boldStart == pStart).
$boldStart = -1
$boldEnd = -1
for ($i = $pStart; $i -lt $pEnd; $i++) {
    $c = $d.Range($i, $i + 1)
    if ($c.Font.Bold -eq -1) {
        if ($boldStart -eq -1) { $boldStart = $i }
        $boldEnd = $i + 1
    }
}

# Replace the bold run's text by inserting the new text just after the old
# bold text and then deleting the old bold text. Doing it in this order (as
# opposed to a direct Range.Text assignment spanning boldStart) avoids
# disturbing the zero-width empty run that immediately precedes the bold
# run at the very start of the paragraph.
$newBoldText = "unfoldingWord® Translation Questions"
$insertPoint = $d.Range($boldEnd, $boldEnd)
$insertPoint.InsertBefore($newBoldText)
$oldBoldRange = $d.Range($boldStart, $boldEnd)
$oldBoldRange.Delete()

$newBoldEnd = $boldStart + $newBoldText.Length

# Everything from the end of the (new) bold run up to, but not including,
# the paragraph mark gets replaced by the new plain-text tail. Use the same
# insert-then-delete trick so the new text doesn't inherit the bold run's
# formatting and so run-splitting stays clean.
$p = $d.Paragraphs.Item($basedIdx)
$pEnd = $p.Range.End
$oldTailEnd = $pEnd - 1

$newTailText = " © 2022 unfoldingWord. Released under CC BY-SA 4.0 license. unfoldingWord® Translation Questions has been adapted in the following languages: Tok Pisin, Arabic (عربي), French (Français), Hindi (हिंदी), Indonesian (Bahasa Indonesia), Portuguese (Português), Russian (Русский), Spanish (Español), Swahili (Kiswahili), and Simplified Chinese (简体中文) from unfoldingWord® Translation Questions © 2022 unfoldingWord. Released under CC BY-SA 4.0 license by Mission Mutual"

$insertPoint2 = $d.Range($oldTailEnd, $oldTailEnd)
$insertPoint2.InsertBefore($newTailText)

$oldTailRange = $d.Range($newBoldEnd, $oldTailEnd)
$oldTailRange.Delete()

Write-Output "Done."
